$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: remove the stale, always-empty "G" placeholder column.
# This shifts H->G, I->H, J->I (values, styles, string refs, and the
# <cols> width entries all move left together), which is exactly the
# structural change in the diff (AAPC/IC_Inf/IC_Sup headers and their
# data move from H:J to G:I).
$ws.Columns("G").Delete()

# Step 2: write the refreshed admissions figures (counts + recomputed
# deltas / AAPC / IC_Inf / IC_Sup) for every data row.
$ws.Range("C2").Value = 5063
$ws.Range("D2").Value = 4393
$ws.Range("E2").Value = 4592
$ws.Range("F2").Value = -9.30278491013233
$ws.Range("G2").Value = -4.87067979136637
$ws.Range("H2").Value = -4.93187072298826
$ws.Range("I2").Value = -4.80944947398879

$ws.Range("C3").Value = 289
$ws.Range("D3").Value = 334
$ws.Range("E3").Value = 427
$ws.Range("F3").Value = 47.7508650519031
$ws.Range("G3").Value = 21.9952739814482
$ws.Range("H3").Value = 21.7055221702211
$ws.Range("I3").Value = 22.2857156226075

$ws.Range("C4").Value = 4774
$ws.Range("D4").Value = 4059
$ws.Range("E4").Value = 4165
$ws.Range("F4").Value = -12.7565982404692
$ws.Range("G4").Value = -6.75549480858085
$ws.Range("H4").Value = -6.8178799443065
$ws.Range("I4").Value = -6.69306790620344

$ws.Range("C5").Value = 4282
$ws.Range("D5").Value = 3774
$ws.Range("E5").Value = 3873
$ws.Range("F5").Value = -9.55161139654367
$ws.Range("G5").Value = -4.99770322224804
$ws.Range("H5").Value = -5.09897336391733
$ws.Range("I5").Value = -4.89632501387013

$ws.Range("C6").Value = 106
$ws.Range("D6").Value = 119
$ws.Range("E6").Value = 177
$ws.Range("F6").Value = 66.9811320754717
$ws.Range("G6").Value = 30.8407909322997
$ws.Range("H6").Value = 30.0673198514908
$ws.Range("I6").Value = 31.6188616120973

$ws.Range("C7").Value = 4177
$ws.Range("D7").Value = 3655
$ws.Range("E7").Value = 3697
$ws.Range("F7").Value = -11.4915010773282
$ws.Range("G7").Value = -6.03780085893659
$ws.Range("H7").Value = -6.13971387826559
$ws.Range("I7").Value = -5.93577718297168

$ws.Range("C8").Value = 4296
$ws.Range("D8").Value = 3655
$ws.Range("E8").Value = 4100
$ws.Range("F8").Value = -4.56238361266294
$ws.Range("G8").Value = -2.33623471196117
$ws.Range("H8").Value = -2.77434675912087
$ws.Range("I8").Value = -1.89614847219419

$ws.Range("C9").Value = 340
$ws.Range("D9").Value = 396
$ws.Range("E9").Value = 564
$ws.Range("F9").Value = 65.8823529411765
$ws.Range("G9").Value = 30.0119294325175
$ws.Range("H9").Value = 28.2123046139125
$ws.Range("I9").Value = 31.8368143031705

$ws.Range("C10").Value = 3956
$ws.Range("D10").Value = 3260
$ws.Range("E10").Value = 3536
$ws.Range("F10").Value = -10.6167846309403
$ws.Range("G10").Value = -5.62582718839986
$ws.Range("H10").Value = -6.07423848823994
$ws.Range("I10").Value = -5.17527512663688

$ws.Range("C11").Value = 6566
$ws.Range("D11").Value = 5588
$ws.Range("E11").Value = 5981
$ws.Range("F11").Value = -8.90953396283887
$ws.Range("G11").Value = -4.60750298742759
$ws.Range("H11").Value = -5.11507427248935
$ws.Range("I11").Value = -4.09721653336038

$ws.Range("C12").Value = 683
$ws.Range("D12").Value = 789
$ws.Range("E12").Value = 956
$ws.Range("F12").Value = 39.9707174231332
$ws.Range("G12").Value = 18.5625923761527
$ws.Range("H12").Value = 16.83257341946
$ws.Range("I12").Value = 20.3182288939667

$ws.Range("C13").Value = 5884
$ws.Range("D13").Value = 4800
$ws.Range("E13").Value = 5025
$ws.Range("F13").Value = -14.5989123045547
$ws.Range("G13").Value = -7.7577701373949
$ws.Range("H13").Value = -8.28548949059822
$ws.Range("I13").Value = -7.22701432121317

$ws.Range("C14").Value = 11086
$ws.Range("D14").Value = 9409
$ws.Range("E14").Value = 9718
$ws.Range("F14").Value = -12.3398881472127
$ws.Range("G14").Value = -6.44527338846238
$ws.Range("H14").Value = -6.93078705917676
$ws.Range("I14").Value = -5.95722694115087

$ws.Range("C15").Value = 1214
$ws.Range("D15").Value = 1424
$ws.Range("E15").Value = 1574
$ws.Range("F15").Value = 29.6540362438221
$ws.Range("G15").Value = 13.6684857799816
$ws.Range("H15").Value = 12.0885859787854
$ws.Range("I15").Value = 15.2706544264846

$ws.Range("C16").Value = 9872
$ws.Range("D16").Value = 7986
$ws.Range("E16").Value = 8144
$ws.Range("F16").Value = -17.5040518638574
$ws.Range("G16").Value = -9.35840909194882
$ws.Range("H16").Value = -9.86591180217631
$ws.Range("I16").Value = -8.84804887235902

# Step 3: restore the cursor/selection the author left the sheet on.
$ws.Range("N22").Select()
